$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly measured / updated performance value for row 18 (Block Inc.)
$ws.Range("C18").Value = -19.89

# Updated "Value (£)" figures recalculated from the refreshed performance data
$ws.Range("D3").Value = 843.6235380817586
$ws.Range("D4").Value = 147.1422559759127
$ws.Range("D5").Value = 157.4233917072693
$ws.Range("D6").Value = 707.2250046894603
$ws.Range("D7").Value = 404.4535414622738
$ws.Range("D8").Value = 237.1946006423465
$ws.Range("D11").Value = 434.393405838844
$ws.Range("D12").Value = 827.426200917405
$ws.Range("D13").Value = 836.9610864709889
$ws.Range("D14").Value = 720.0001932963053
$ws.Range("D15").Value = 436.9474943508601
$ws.Range("D16").Value = 542.5287923954203
$ws.Range("D17").Value = 68.35959057087283
$ws.Range("D18").Value = 8010.756556715879
